# NYPD CompStat weekly report refresh: new crime data collected for the
# week of 4/17/2023 - 4/23/2023 (previously 4/10/2023 - 4/16/2023), Volume 30
# Number 16 (previously Number 15). Updates the report header text and the
# full Bronx crime-complaint statistics table (rows 14-30: Week to Date,
# 28 Day, Year to Date, and 2 Year figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: issue number and covered week ---
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# --- Crime complaint statistics table (rows 14-30) ---

# Row 14: Murder
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0
$ws.Range("I14").Value = 35
$ws.Range("J14").Value = 45
$ws.Range("K14").Value = -22.222222222222
$ws.Range("L14").Value = -10.256410256410
$ws.Range("M14").Value = 2.941176470588
$ws.Range("N14").Value = -76.666666666666

# Row 15: Rape
$ws.Range("C15").Value = 11
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 83.333333333333
$ws.Range("F15").Value = 34
$ws.Range("G15").Value = 27
$ws.Range("H15").Value = 25.925925925925
$ws.Range("I15").Value = 125
$ws.Range("J15").Value = 124
$ws.Range("K15").Value = 0.806451612903
$ws.Range("L15").Value = 20.192307692307
$ws.Range("M15").Value = 52.439024390243
$ws.Range("N15").Value = -38.118811881188

# Row 16: Robbery
$ws.Range("C16").Value = 74
$ws.Range("D16").Value = 86
$ws.Range("E16").Value = -13.953488372093
$ws.Range("F16").Value = 316
$ws.Range("G16").Value = 315
$ws.Range("H16").Value = 0.317460317460
$ws.Range("I16").Value = 1310
$ws.Range("J16").Value = 1318
$ws.Range("K16").Value = -0.606980273141
$ws.Range("L16").Value = 40.708915145005
$ws.Range("M16").Value = 6.764466177669
$ws.Range("N16").Value = -73.820943245403

# Row 17: Fel. Assault
$ws.Range("C17").Value = 137
$ws.Range("D17").Value = 122
$ws.Range("E17").Value = 12.295081967213
$ws.Range("F17").Value = 637
$ws.Range("G17").Value = 538
$ws.Range("H17").Value = 18.401486988847
$ws.Range("I17").Value = 2204
$ws.Range("J17").Value = 1986
$ws.Range("K17").Value = 10.976837865055
$ws.Range("L17").Value = 36.724565756823
$ws.Range("M17").Value = 71.784879189399
$ws.Range("N17").Value = -7.820995399414

# Row 18: Burglary
$ws.Range("C18").Value = 58
$ws.Range("D18").Value = 48
$ws.Range("E18").Value = 20.833333333333
$ws.Range("F18").Value = 244
$ws.Range("G18").Value = 224
$ws.Range("H18").Value = 8.928571428571
$ws.Range("I18").Value = 955
$ws.Range("J18").Value = 923
$ws.Range("K18").Value = 3.466955579631
$ws.Range("L18").Value = 47.376543209876
$ws.Range("M18").Value = 4.143947655398
$ws.Range("N18").Value = -83.082373782108

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 168
$ws.Range("D19").Value = 141
$ws.Range("E19").Value = 19.148936170212
$ws.Range("F19").Value = 632
$ws.Range("G19").Value = 605
$ws.Range("H19").Value = 4.462809917355
$ws.Range("I19").Value = 2305
$ws.Range("J19").Value = 2416
$ws.Range("K19").Value = -4.594370860927
$ws.Range("L19").Value = 32.091690544412
$ws.Range("M19").Value = 83.665338645418
$ws.Range("N19").Value = 11.676356589147

# Row 20: G.L.A.
$ws.Range("C20").Value = 96
$ws.Range("D20").Value = 77
$ws.Range("E20").Value = 24.675324675324
$ws.Range("F20").Value = 403
$ws.Range("G20").Value = 268
$ws.Range("H20").Value = 50.373134328358
$ws.Range("I20").Value = 1614
$ws.Range("J20").Value = 1300
$ws.Range("K20").Value = 24.153846153846
$ws.Range("L20").Value = 129.587482219061
$ws.Range("M20").Value = 179.722703639515
$ws.Range("N20").Value = -66.149328859060

# Row 21: TOTAL
$ws.Range("C21").Value = 547
$ws.Range("D21").Value = 483
$ws.Range("E21").Value = 13.250517598343
$ws.Range("F21").Value = 2274
$ws.Range("G21").Value = 1988
$ws.Range("H21").Value = 14.386317907444
$ws.Range("I21").Value = 8548
$ws.Range("J21").Value = 8112
$ws.Range("K21").Value = 5.374753451676
$ws.Range("L21").Value = 47.838118298166
$ws.Range("M21").Value = 59.032558139534
$ws.Range("N21").Value = -57.733386075949

# Row 22: Transit
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 9
$ws.Range("E22").Value = -11.111111111111
$ws.Range("F22").Value = 29
$ws.Range("G22").Value = 24
$ws.Range("H22").Value = 20.833333333333
$ws.Range("I22").Value = 88
$ws.Range("J22").Value = 103
$ws.Range("K22").Value = -14.563106796116
$ws.Range("L22").Value = 29.411764705882
$ws.Range("M22").Value = -8.333333333333

# Row 23: Housing
$ws.Range("C23").Value = 37
$ws.Range("D23").Value = 36
$ws.Range("E23").Value = 2.777777777777
$ws.Range("F23").Value = 158
$ws.Range("G23").Value = 122
$ws.Range("H23").Value = 29.508196721311
$ws.Range("I23").Value = 544
$ws.Range("J23").Value = 451
$ws.Range("K23").Value = 20.620842572062
$ws.Range("L23").Value = 66.871165644171
$ws.Range("M23").Value = 88.888888888888

# Row 24: Petit Larceny
$ws.Range("C24").Value = 335
$ws.Range("D24").Value = 344
$ws.Range("E24").Value = -2.616279069767
$ws.Range("F24").Value = 1408
$ws.Range("G24").Value = 1483
$ws.Range("H24").Value = -5.057316250842
$ws.Range("I24").Value = 5292
$ws.Range("J24").Value = 5315
$ws.Range("K24").Value = -0.432737535277
$ws.Range("L24").Value = 44.629680240502
$ws.Range("M24").Value = 46.877601998334

# Row 25: Misd. Assault
$ws.Range("C25").Value = 196
$ws.Range("D25").Value = 173
$ws.Range("E25").Value = 13.294797687861
$ws.Range("F25").Value = 811
$ws.Range("G25").Value = 767
$ws.Range("H25").Value = 5.736636245110
$ws.Range("I25").Value = 3050
$ws.Range("J25").Value = 2911
$ws.Range("K25").Value = 4.774991411885
$ws.Range("L25").Value = 31.408875484704
$ws.Range("M25").Value = -1.739690721649

# Row 26: UCR Rape*
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 128.571428571429
$ws.Range("F26").Value = 52
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = 8.333333333333
$ws.Range("I26").Value = 196
$ws.Range("J26").Value = 218
$ws.Range("K26").Value = -10.091743119266
$ws.Range("L26").Value = 11.363636363636

# Row 27: Other Sex Crimes
$ws.Range("C27").Value = 24
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 140
$ws.Range("F27").Value = 83
$ws.Range("G27").Value = 70
$ws.Range("H27").Value = 18.571428571428
$ws.Range("I27").Value = 323
$ws.Range("J27").Value = 266
$ws.Range("K27").Value = 21.428571428571
$ws.Range("L27").Value = 19.188191881918

# Row 28: Shooting Vic.
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = -58.333333333333
$ws.Range("F28").Value = 27
$ws.Range("G28").Value = 46
$ws.Range("H28").Value = -41.304347826087
$ws.Range("I28").Value = 99
$ws.Range("J28").Value = 154
$ws.Range("K28").Value = -35.714285714285
$ws.Range("L28").Value = -14.655172413793
$ws.Range("M28").Value = -13.913043478260
$ws.Range("N28").Value = -71.633237822349

# Row 29: Shooting Inc.
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 8
$ws.Range("E29").Value = -37.5
$ws.Range("F29").Value = 21
$ws.Range("G29").Value = 33
$ws.Range("H29").Value = -36.363636363636
$ws.Range("I29").Value = 79
$ws.Range("J29").Value = 132
$ws.Range("K29").Value = -40.151515151515
$ws.Range("L29").Value = -26.168224299065
$ws.Range("M29").Value = -21
$ws.Range("N29").Value = -75

# Row 30: Hate Crimes
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 9
$ws.Range("K30").Value = -40
$ws.Range("L30").Value = -18.181818181818
